$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source report had several adjacent row-pairs whose B:G contents
# (Code, Name, MRP, Price, Qty, Amount) were transposed between the two
# rows. Column A (Sr No) is correct and stays untouched; for each pair we
# swap the B:G block between row1 and row2.
$rows1 = @(313, 316, 346, 351, 355, 379, 382, 389, 400, 419, 421, 431, 457, 536, 579, 583, 586, 593, 599, 604, 720, 872)
$rows2 = @(314, 317, 347, 352, 356, 380, 383, 390, 401, 420, 422, 432, 458, 537, 580, 584, 587, 594, 600, 605, 721, 873)

for ($i = 0; $i -lt $rows1.Length; $i++) {
    $row1 = $rows1[$i]
    $row2 = $rows2[$i]

    $range1 = $ws.Range("B$row1`:G$row1")
    $range2 = $ws.Range("B$row2`:G$row2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
